{"js": "// Add a two-space indentation run in front of every table-cell paragraph\n// whose text looks like a \"feature : Type\" line (e.g. \"companies : Company\"),\n// mirroring the Sirius table conversion's new sub-line indentation support.\n// Class-header rows (\"World\"), inheritance rows (\"Company -> NamedElement\")\n// and plain type names are left untouched because their text has no \" : \".\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\n// Collect every paragraph whose text contains the \"name : Type\" pattern.\nconst targets = [];\nfor (let i = 0; i < body.paragraphs.items.length; i++) {\n  const paragraph = body.paragraphs.items[i];\n  if (paragraph.text.indexOf(\" : \") !== -1) {\n    targets.push(paragraph);\n  }\n}\n\n// Insert a \"  \" run at the start of each target paragraph, copying the\n// bold/size/color formatting of the paragraph's existing run. Flip italic\n// on then off around the sync boundary so the new run is kept distinct\n// from the following run instead of being silently re-merged with it\n// (both runs end up with identical formatting, exactly like the source).\nconst inserted = [];\nfor (const paragraph of targets) {\n  const spacer = paragraph.insertText(\"  \", \"Start\");\n  spacer.font.bold = true;\n  spacer.font.size = 12;\n  spacer.font.color = \"#000000\";\n  spacer.font.italic = true;\n  inserted.push(spacer);\n}\nawait context.sync();\n\nfor (const spacer of inserted) {\n  spacer.font.italic = false;\n}\nawait context.sync();\n", "ps1": "# Add a two-space indentation run in front of every table-cell paragraph\n# whose text looks like a \"feature : Type\" line (e.g. \"companies : Company\"),\n# mirroring the Sirius table conversion's new sub-line indentation support.\n# Class-header rows (\"World\"), inheritance rows (\"Company -> NamedElement\")\n# and plain type names are left untouched because their text has no \" : \".\n\n$d = $word.ActiveDocument\n\n# Collect the start offset of every paragraph whose text contains the\n# \"name : Type\" pattern, before any edits shift document positions.\n$targets = @()\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"* : *\") {\n        $targets += $p.Range.Start\n    }\n}\n\n# Walk the offsets back-to-front so inserting earlier in the document never\n# invalidates an offset we still need to use.\nfor ($j = $targets.Count - 1; $j -ge 0; $j--) {\n    $startPos = $targets[$j]\n\n    $r = $d.Range($startPos, $startPos)\n    $r.InsertBefore(\"  \")\n\n    # Copy the bold/size/color formatting of the line onto the new spacer\n    # run. Flip Italic on then back off so the engine keeps the spacer as\n    # its own run instead of silently re-merging it with the run that\n    # follows (both end up with identical formatting, same as the source).\n    $ins = $d.Range($startPos, $startPos + 2)\n    $ins.Font.Bold = 1\n    $ins.Font.Size = 12\n    $ins.Font.Color = 0\n    $ins.Font.Italic = 1\n    $ins.Font.Italic = 0\n}\n"}
